$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header: new "Score Set 4" column (H)
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "Score Set 4"

# ---------------------------------------------------------------------------
# Row 8 - Elizabethtown & Lebanon & Northern Lebanon @ Hempfield (Boys)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Elizabethtown & Lebanon & Northern Lebanon @ Hempfield (Boys)"
$ws.Range("A9").Value = "Elizabethtown & Lebanon & Northern Lebanon @ Hempfield (Girls)"

$ws.Range("C8").Value = "Mid-70s & Clear Skies"
$ws.Range("C9").Value = "Mid-70s & Clear Skies"

$ws.Range("D8").Value = "We had a good race, but we need to step up our game & close that 1:46 gap between #1 & #5 runners in the upcoming races.  Personally, my goal is to have a 60 second spread 1-5 off of Galli & a 20 second gap between #2 & #5.  We can do it...let's get after it!"
$ws.Range("D9").Value = "I loved the effort!  Let's get better on Saturday!"

$ws.Range("E8").Value = "Cedar Crest 23 - Hempfield 34"
$ws.Range("F8").Value = "Cedar Crest 21 - Elizabethtown 37"
$ws.Range("G8").Value = "Cedar Crest 15 - Lebanon 55"
$ws.Range("H8").Value = "Cedar Crest 15 - Northern Lebanon 55"

$ws.Range("E9").Value = "Cedar Crest 36 - Hempfield 19"
$ws.Range("F9").Value = "Cedar Crest 32 - Elizabethtown 24"
$ws.Range("G9").Value = "Cedar Crest 15 - Lebanon 55"
$ws.Range("H9").Value = "Cedar Crest 15 - Northern Lebanon 55"

# ---------------------------------------------------------------------------
# Rows 10-13 - Lebanon County Meet @ South Hills
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "Lebanon County Meet @ South Hills (Varsity Boys)"
$ws.Range("A11").Value = "Lebanon County Meet @ South Hills (JV Boys)"
$ws.Range("A12").Value = "Lebanon County Meet @ South Hills (Varsity Girls)"
$ws.Range("A13").Value = "Lebanon County Meet @ South Hills (JV Girls)"

$ws.Range("C10").Value = "Mid-60s, Sunny and a Little Breezy"
$ws.Range("C11").Value = "Mid-60s, Sunny and a Little Breezy"
$ws.Range("C12").Value = "Mid-60s, Sunny and a Little Breezy"
$ws.Range("C13").Value = "Mid-60s, Sunny and a Little Breezy"

$ws.Range("D10").Value = "Team Results: FIRST PLACE, 28 Points. #1-#5 Spread= 1:20…getting better, but we need to get down to 1:00 or under!"
$ws.Range("D11").Value = "Team Results: 1st-9th!"
$ws.Range("D12").Value = "TEAM Results: 2nd Place, 50 Points"
$ws.Range("D13").Value = "?"

# ---------------------------------------------------------------------------
# Rows 14-15 - Conestoga Valley & Warkwick & Garden Spot & Ephrata @ Home
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Conestoga Valley & Warkwick & Garden Spot & Ephrata @ Home (Boys)"
$ws.Range("A15").Value = "Conestoga Valley & Warkwick & Garden Spot & Ephrata @ Home (Girls)"

$ws.Range("C14").Value = "?"
$ws.Range("C15").Value = "?"
$ws.Range("D14").Value = "?"
$ws.Range("D15").Value = "?"

# ---------------------------------------------------------------------------
# Rows 16-19 - Carlisle Meet of Champions
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "Carlisle Meet of Champions (Varsity Boys)"
$ws.Range("A17").Value = "Carlisle Meet of Champions (JV Boys)"
$ws.Range("A18").Value = "Carlisle Meet of Champions (Varsity Girls)"
$ws.Range("A19").Value = "Carlisle Meet of Champions (JV Girls)"

$ws.Range("C16").Value = "?"
$ws.Range("C17").Value = "?"
$ws.Range("C18").Value = "?"
$ws.Range("C19").Value = "?"
$ws.Range("D16").Value = "?"
$ws.Range("D17").Value = "?"
$ws.Range("D18").Value = "?"
$ws.Range("D19").Value = "?"

# ---------------------------------------------------------------------------
# Rows 20-21 - Lancaster Mennonite & Columbia & Lampeter-Strasburg @ Donegal
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "Lancaster Mennonite & Columbia & Lampeter-Strasburg @ Donegal (Boys)"
$ws.Range("A21").Value = "Lancaster Mennonite & Columbia & Lampeter-Strasburg @ Donegal (Girls)"

$ws.Range("C20").Value = "?"
$ws.Range("C21").Value = "?"
$ws.Range("D20").Value = "?"
$ws.Range("D21").Value = "?"

$ws.Range("E20").Value = "Cedar Crest 15 - Donegal 55"
$ws.Range("F20").Value = "Cedar Crest 15 - Columbia 55"
$ws.Range("G20").Value = "Cedar Crest 23 - Lancaster Mennonite 34"
$ws.Range("H20").Value = "Cedar Crest 17 - Lampeter-Strasburg 44"

$ws.Range("E21").Value = "Cedar Crest 15 - Donegal 55"
$ws.Range("F21").Value = "Cedar Crest 15 - Columbia 55"
$ws.Range("G21").Value = "Cedar Crest 15 - Lancaster Mennonite 55"
$ws.Range("H21").Value = "Cedar Crest 17 - Lampeter-Strasburg 44"

# ---------------------------------------------------------------------------
# Dates (column B) - copy number format from an existing date cell (B2) so
# the style index matches instead of minting a brand-new number format.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B8:B9").PasteSpecial(-4122)
$ws.Range("B10:B13").PasteSpecial(-4122)
$ws.Range("B14:B15").PasteSpecial(-4122)
$ws.Range("B16:B19").PasteSpecial(-4122)
$ws.Range("B20:B21").PasteSpecial(-4122)

$ws.Range("B8").Value = 40435
$ws.Range("B9").Value = 40435
$ws.Range("B10").Value = 40439
$ws.Range("B11").Value = 40439
$ws.Range("B12").Value = 40439
$ws.Range("B13").Value = 40439
$ws.Range("B14").Value = 40442
$ws.Range("B15").Value = 40442
$ws.Range("B16").Value = 40446
$ws.Range("B17").Value = 40446
$ws.Range("B18").Value = 40446
$ws.Range("B19").Value = 40446
$ws.Range("B20").Value = 40449
$ws.Range("B21").Value = 40449

# ---------------------------------------------------------------------------
# Column D (Coach's Comments) font style - copy style from D2 (s="1") for
# the rows that have it; D10/D12 keep the default style (no explicit s).
# ---------------------------------------------------------------------------
$ws.Range("D2").Copy()
$ws.Range("D8:D9").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D13:D19").PasteSpecial(-4122)
$ws.Range("D20:D21").PasteSpecial(-4122)

$ws.Range("D8").Value = "We had a good race, but we need to step up our game & close that 1:46 gap between #1 & #5 runners in the upcoming races.  Personally, my goal is to have a 60 second spread 1-5 off of Galli & a 20 second gap between #2 & #5.  We can do it...let's get after it!"
$ws.Range("D9").Value = "I loved the effort!  Let's get better on Saturday!"
$ws.Range("D11").Value = "Team Results: 1st-9th!"
$ws.Range("D13").Value = "?"
$ws.Range("D14").Value = "?"
$ws.Range("D15").Value = "?"
$ws.Range("D16").Value = "?"
$ws.Range("D17").Value = "?"
$ws.Range("D18").Value = "?"
$ws.Range("D19").Value = "?"
$ws.Range("D21").Value = "?"

# ---------------------------------------------------------------------------
# Column widths (COM ColumnWidth is quantized to ~1/6-character steps by the
# host, so these are the closest achievable values to the target XML widths
# 27.5703125 / 30.7109375 / 33.42578125)
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 26.6
$ws.Columns("F").ColumnWidth = 29.75
$ws.Columns("H").ColumnWidth = 32.6

# ---------------------------------------------------------------------------
# Selection
# ---------------------------------------------------------------------------
$ws.Range("D10").Select()
